$d = $word.ActiveDocument

# Replaces the first occurrence of $target that is found inside the (unique)
# surrounding $context string, by locating $context in the live document
# text, then locating $target inside it, and finally writing the new text
# into a precise Range that covers only $target. Using extra (unmodified)
# context characters keeps the match anchored to the correct occurrence
# without letting the edited span bleed into neighbouring runs.
function Replace-Span($context, $target, $newText) {
    $full = $d.Content.Text
    $ctxIdx = $full.IndexOf($context)
    if ($ctxIdx -lt 0) {
        throw "context not found: $context"
    }
    $targetIdx = $context.IndexOf($target)
    if ($targetIdx -lt 0) {
        throw "target not found inside context: $target"
    }
    $start = $ctxIdx + $targetIdx
    $end = $start + $target.Length
    $r = $d.Range($start, $end)
    $r.Text = $newText
}

# "Decompomos ... como por exemplo Employee, ..." -> "... como por exemplo a de Employee, ..."
Replace-Span "como por exemplo Employee" "como por exemplo " "como por exemplo a de "

# "Director, analyst, manager, especialist." -> "Director, Manager, Especialist e Analyst."
Replace-Span "Director, analyst" ", " ", Manager, "
Replace-Span "analyst" "analyst" "Especialist"
Replace-Span ", manager, " ", manager, " " e "
Replace-Span "especialist" "especialist" "Analyst"

# "... que compõem a entidade employee preservam ..." -> "... entidade Employee preservam ..."
Replace-Span "entidade employee preservam" "employee" "Employee"

# "No começo, a entidade TEAM tinha ..." -> "... a entidade Team tinha ..."
Replace-Span "a entidade TEAM tinha" "TEAM" "Team"

# "... removemos ele para deixar na terceira forma normal. " ->
# "... removemos ele para deixar na terceira forma normal assim como todos
#     os outros atributos que tinham a mesma finalidade presentes em outras
#     entidades. "
Replace-Span "terceira forma normal. " "terceira forma normal. " `
    "terceira forma normal assim como todos os outros atributos que tinham a mesma finalidade presentes em outras entidades. "
